$wb = $excel.ActiveWorkbook

$wsP = $wb.Worksheets.Item("P_valores")
$wsP.Range("B3").Value = 0.03928680844012988
$wsP.Range("B4").Value = 0.9207881778545348
$wsP.Range("B5").Value = 0.5237534330737024
$wsP.Range("B6").Value = 0.08934961645053097
$wsP.Range("C2").Value = 0.03928680844012988
$wsP.Range("C4").Value = 0.01521214829376549
$wsP.Range("C5").Value = 0.1768162151734904
$wsP.Range("C6").Value = 0.4714513532438007
$wsP.Range("D2").Value = 0.9207881778545348
$wsP.Range("D3").Value = 0.01521214829376549
$wsP.Range("D5").Value = 0.5927375633370369
$wsP.Range("D6").Value = 0.03005681517609471
$wsP.Range("E2").Value = 0.5237534330737024
$wsP.Range("E3").Value = 0.1768162151734904
$wsP.Range("E4").Value = 0.5927375633370369
$wsP.Range("E6").Value = 0.1967716791863061
$wsP.Range("F2").Value = 0.08934961645053097
$wsP.Range("F3").Value = 0.4714513532438007
$wsP.Range("F4").Value = 0.03005681517609471
$wsP.Range("F5").Value = 0.1967716791863061

$wsE = $wb.Worksheets.Item("Estadisticos_DM")
$wsE.Range("B3").Value = 2.273302295993683
$wsE.Range("B4").Value = -0.1012495069954851
$wsE.Range("B5").Value = 0.6539304745123241
$wsE.Range("B6").Value = 1.825368721786241
$wsE.Range("C2").Value = -2.273302295993683
$wsE.Range("C4").Value = -2.764283591989347
$wsE.Range("C5").Value = -1.422369511159142
$wsE.Range("C6").Value = 0.7401305020702474
$wsE.Range("D2").Value = 0.1012495069954851
$wsE.Range("D3").Value = 2.764283591989347
$wsE.Range("D5").Value = 0.5473795579644489
$wsE.Range("D6").Value = 2.413911161704107
$wsE.Range("E2").Value = -0.6539304745123241
$wsE.Range("E3").Value = 1.422369511159142
$wsE.Range("E4").Value = -0.5473795579644489
$wsE.Range("E6").Value = 1.355345895474942
$wsE.Range("F2").Value = -1.825368721786241
$wsE.Range("F3").Value = -0.7401305020702474
$wsE.Range("F4").Value = -2.413911161704107
$wsE.Range("F5").Value = -1.355345895474942
